$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '51.206.07'
$ws.Range("E2").Value = '  +2.97%  '

# Row 3
$ws.Range("D3").Value = '2.744.97'
$ws.Range("E3").Value = '  +3.01%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '115.13'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.81%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '333.38'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.50%  '

# Row 7
$ws.Range("E7").Value = '  +0.66%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.09%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.572'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.63%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.36'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.06%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.14'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.49%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0828'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.79%  '

# Row 13
$ws.Range("E13").Value = '  +2.96%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.65'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.23%  '

# Row 15
$ws.Range("D15").Value = '3.173.47'
$ws.Range("E15").Value = '  +2.69%  '

# Row 16
$ws.Range("D16").Value = '2.733.70'
$ws.Range("E16").Value = '  +1.98%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.885'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.52%  '

# Row 18
$ws.Range("D18").Value = '51.083.64'
$ws.Range("E18").Value = '  +2.74%  '

# Row 19
$ws.Range("E19").Value = '  +3.84%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.01'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.22%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.86'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.32%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0960'
$ws.Range("E22").Value = '  +0.51%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '279.60'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.18%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.17'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.64%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.63'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.55%  '

# Row 26
$ws.Range("E26").Value = '  +0.73%  '

# Row 27
$ws.Range("E27").Value = '  +0.20%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.35'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.51%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.56%  '

# Row 30
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.63'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.16%  '

# Row 31
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.140'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.23%  '

# Row 32
$ws.Range("E32").Value = '  -0.46%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.61'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.45%  '

# Row 34
$ws.Range("E34").Value = '  +1.87%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.38'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.63%  '

# Row 36
$ws.Range("E36").Value = '  -0.34%  '

# Row 37
$ws.Range("E37").Value = '  +2.13%  '

# Row 38
$ws.Range("E38").Value = '  -0.37%  '

# Row 39
$ws.Range("E39").Value = '  +2.05%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '129.30'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.90%  '

# Row 41
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.69'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +4.39%  '

# Row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0353'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +11.21%  '

# Row 43
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.114'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.86%  '

# Row 44
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.28'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.83%  '

# Row 45
$ws.Range("D45").Value = '2.113.74'
$ws.Range("E45").Value = '  +0.53%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.39'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.13%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.23'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +9.36%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.30'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.63%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.54'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.51%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.07'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.81%  '

# Row 51
$ws.Range("E51").Value = '  +10.49%  '
